# Auto-assembled Excel COM-interop script
# Restructures the single-sheet "ODI Batting" workbook into three sheets:
#   "Player Info", "ODI Batting" (updated), "ODI Batting Extra"
$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $val) {
    # Writes $val as a cell value. Callers pass numeric-looking
    # text already prefixed with a leading apostrophe so Excel
    # keeps it as text instead of converting it to a number.
    if ($null -eq $val) { return }
    $ws.Cells.Item($row, $col).Value = $val
}

function Set-NumCell($ws, $row, $col, $val) {
    if ($null -eq $val) { return }
    $ws.Cells.Item($row, $col).Value = $val
}

function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
    $range.Borders.LineStyle = 1
}

# --- Update existing "ODI Batting" sheet: header + MATCH_CODE column, drop stray empty cells ---
$odiBatting = $wb.Worksheets.Item(1)
Set-TextCell $odiBatting 1 4 "MATCH_CODE"

Set-TextCell $odiBatting 2 4 "'3162"
Set-TextCell $odiBatting 3 4 "'3168"
Set-TextCell $odiBatting 4 4 "'3170"
Set-TextCell $odiBatting 5 4 "'3182"
Set-TextCell $odiBatting 6 4 "'3186"
Set-TextCell $odiBatting 7 4 "'3188"
Set-TextCell $odiBatting 8 4 "'3190"
Set-TextCell $odiBatting 9 4 "'3342"
Set-TextCell $odiBatting 10 4 "'3345"
Set-TextCell $odiBatting 11 4 "'3422"
Set-TextCell $odiBatting 12 4 "'3424"
Set-TextCell $odiBatting 13 4 "'3427"
Set-TextCell $odiBatting 14 4 "'3446"
Set-TextCell $odiBatting 15 4 "'3447"
Set-TextCell $odiBatting 16 4 "'3448"
Set-TextCell $odiBatting 17 4 "'3449"
Set-TextCell $odiBatting 18 4 "'3464"
Set-TextCell $odiBatting 19 4 "'3466"
Set-TextCell $odiBatting 20 4 "'3469"
Set-TextCell $odiBatting 21 4 "'3476"
Set-TextCell $odiBatting 22 4 "'3477"
Set-TextCell $odiBatting 23 4 "'3479"
Set-TextCell $odiBatting 24 4 "'3954"
Set-TextCell $odiBatting 25 4 "'3955"
Set-TextCell $odiBatting 26 4 "'3966"
Set-TextCell $odiBatting 27 4 "'3967"
Set-TextCell $odiBatting 28 4 "'3968"
Set-TextCell $odiBatting 29 4 "'4227"

# Remove stray empty INNING_NUMBER cells (rows where inning number is blank)
$odiBatting.Cells.Item(3, 2).ClearContents()
$odiBatting.Cells.Item(17, 2).ClearContents()
$odiBatting.Cells.Item(29, 2).ClearContents()

# --- Add "Player Info" sheet before "ODI Batting" ---
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

Set-TextCell $playerInfo 1 1 "ID"
Set-TextCell $playerInfo 1 2 "NAME"
Set-TextCell $playerInfo 1 3 "BATTING_HAND"
Set-TextCell $playerInfo 1 4 "BOWL_STYLE"
Set-HeaderStyle $playerInfo.Range("A1:D1")

Set-TextCell $playerInfo 2 1 "'3737"
Set-TextCell $playerInfo 2 2 "Bradley-John Watling"
Set-TextCell $playerInfo 2 3 "Right Handed"
Set-TextCell $playerInfo 2 4 "Does Not Bowl | Unknown"

# --- Add "ODI Batting Extra" sheet after "ODI Batting" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

Set-TextCell $extra 1 1 "MATCH_CODE"
Set-TextCell $extra 1 2 "BATTING_POSITION"
Set-TextCell $extra 1 3 "NUM_4"
Set-TextCell $extra 1 4 "NUM_6"
Set-TextCell $extra 1 5 "PERCENT_RUNS_OF_TOTAL"
Set-TextCell $extra 1 6 "MAN_OF_MATCH"
Set-HeaderStyle $extra.Range("A1:F1")

Set-TextCell $extra 2 1 "'3345"
Set-NumCell $extra 2 2 2
Set-TextCell $extra 2 3 "'0"
Set-TextCell $extra 2 4 "'0"
Set-TextCell $extra 2 5 $null
Set-TextCell $extra 2 6 "NO"

Set-TextCell $extra 3 1 "'3422"
Set-NumCell $extra 3 2 6
Set-TextCell $extra 3 3 "'5"
Set-TextCell $extra 3 4 "'0"
Set-TextCell $extra 3 5 "'31.58%"
Set-TextCell $extra 3 6 "NO"

Set-TextCell $extra 4 1 "'3424"
Set-NumCell $extra 4 2 6
Set-TextCell $extra 4 3 "'4"
Set-TextCell $extra 4 4 "'2"
Set-TextCell $extra 4 5 "'27.69%"
Set-TextCell $extra 4 6 "NO"

Set-TextCell $extra 5 1 "'3427"
Set-NumCell $extra 5 2 6
Set-TextCell $extra 5 3 "'4"
Set-TextCell $extra 5 4 "'1"
Set-TextCell $extra 5 5 "'16.06%"
Set-TextCell $extra 5 6 "NO"

Set-TextCell $extra 6 1 "'3446"
Set-NumCell $extra 6 2 3
Set-TextCell $extra 6 3 "'6"
Set-TextCell $extra 6 4 "'0"
Set-TextCell $extra 6 5 "'22.00%"
Set-TextCell $extra 6 6 "NO"

Set-TextCell $extra 7 1 "'3447"
Set-NumCell $extra 7 2 2
Set-TextCell $extra 7 3 "'12"
Set-TextCell $extra 7 4 "'0"
Set-TextCell $extra 7 5 "'51.06%"
Set-TextCell $extra 7 6 "NO"

Set-TextCell $extra 8 1 "'3448"
Set-TextCell $extra 8 3 $null
Set-TextCell $extra 8 4 $null
Set-TextCell $extra 8 5 $null
Set-TextCell $extra 8 6 "NO"

Set-TextCell $extra 9 1 "'3449"
Set-TextCell $extra 9 3 $null
Set-TextCell $extra 9 4 $null
Set-TextCell $extra 9 5 $null
Set-TextCell $extra 9 6 "NO"

Set-TextCell $extra 10 1 "'3464"
Set-NumCell $extra 10 2 3
Set-TextCell $extra 10 3 "'5"
Set-TextCell $extra 10 4 "'0"
Set-TextCell $extra 10 5 "'21.53%"
Set-TextCell $extra 10 6 "NO"

Set-TextCell $extra 11 1 "'3466"
Set-TextCell $extra 11 3 $null
Set-TextCell $extra 11 4 $null
Set-TextCell $extra 11 5 $null
Set-TextCell $extra 11 6 "NO"

Set-TextCell $extra 12 1 "'3469"
Set-NumCell $extra 12 2 2
Set-TextCell $extra 12 3 "'2"
Set-TextCell $extra 12 4 "'0"
Set-TextCell $extra 12 5 "'7.69%"
Set-TextCell $extra 12 6 "NO"

Set-TextCell $extra 13 1 "'3476"
Set-NumCell $extra 13 2 1
Set-TextCell $extra 13 3 "'0"
Set-TextCell $extra 13 4 "'0"
Set-TextCell $extra 13 5 "'0.77%"
Set-TextCell $extra 13 6 "NO"

Set-TextCell $extra 14 1 "'3477"
Set-NumCell $extra 14 2 1
Set-TextCell $extra 14 3 "'0"
Set-TextCell $extra 14 4 "'0"
Set-TextCell $extra 14 5 "'2.60%"
Set-TextCell $extra 14 6 "NO"

Set-TextCell $extra 15 1 "'3479"
Set-NumCell $extra 15 2 1
Set-TextCell $extra 15 3 "'0"
Set-TextCell $extra 15 4 "'0"
Set-TextCell $extra 15 5 "'0.54%"
Set-TextCell $extra 15 6 "NO"

Set-TextCell $extra 16 1 "'3954"
Set-NumCell $extra 16 2 6
Set-TextCell $extra 16 3 "'1"
Set-TextCell $extra 16 4 "'0"
Set-TextCell $extra 16 5 "'5.38%"
Set-TextCell $extra 16 6 "NO"

Set-TextCell $extra 17 1 "'3955"
Set-NumCell $extra 17 2 6
Set-TextCell $extra 17 3 "'0"
Set-TextCell $extra 17 4 "'0"
Set-TextCell $extra 17 5 $null
Set-TextCell $extra 17 6 "NO"

Set-TextCell $extra 18 1 "'3966"
Set-TextCell $extra 18 3 $null
Set-TextCell $extra 18 4 $null
Set-TextCell $extra 18 5 $null
Set-TextCell $extra 18 6 "NO"

Set-TextCell $extra 19 1 "'3967"
Set-NumCell $extra 19 2 8
Set-TextCell $extra 19 3 "'2"
Set-TextCell $extra 19 4 "'0"
Set-TextCell $extra 19 5 "'6.49%"
Set-TextCell $extra 19 6 "NO"

Set-TextCell $extra 20 1 "'3968"
Set-TextCell $extra 20 3 $null
Set-TextCell $extra 20 4 $null
Set-TextCell $extra 20 5 $null
Set-TextCell $extra 20 6 "NO"

Set-TextCell $extra 21 1 "'4227"
Set-TextCell $extra 21 3 $null
Set-TextCell $extra 21 4 $null
Set-TextCell $extra 21 5 $null
Set-TextCell $extra 21 6 "NO"

Write-Host "Edit complete"
